# Scheduled market-data refresh: update cached price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across all job sheets, per latest Universalis pull.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 255.1579
$ws.Range("I33").Value = 244.58824
$ws.Range("J33").Value = 345
$ws.Range("K33").Value = 244.58824
$ws.Range("L33").Value = 345
$ws.Range("M33").Value = -15.58824000000001
$ws.Range("N33").Value = -803
$ws.Range("H62").Value = 4467.25
$ws.Range("I62").Value = 3029.3333
$ws.Range("J62").Value = 5905.1665
$ws.Range("K62").Value = 3029.3333
$ws.Range("L62").Value = 5905.1665
$ws.Range("M62").Value = -2405.3333
$ws.Range("N62").Value = -7153.1665
$ws.Range("H65").Value = 4467.25
$ws.Range("I65").Value = 3029.3333
$ws.Range("J65").Value = 5905.1665
$ws.Range("K65").Value = 15146.6665
$ws.Range("L65").Value = 29525.8325
$ws.Range("M65").Value = -12026.6665
$ws.Range("N65").Value = -35765.8325
$ws.Range("H112").Value = 1319.7354
$ws.Range("J112").Value = 1370.9688
$ws.Range("L112").Value = 4112.9064
$ws.Range("N112").Value = -6328.9064
$ws.Range("H129").Value = 847.76404
$ws.Range("J129").Value = 904.65
$ws.Range("L129").Value = 2713.95
$ws.Range("N129").Value = -12713.95
$ws.Range("H132").Value = 34832828
$ws.Range("I132").Value = 38466724
$ws.Range("K132").Value = 115400172
$ws.Range("M132").Value = -115397642
$ws.Range("H137").Value = 1445191
$ws.Range("I137").Value = 2507851.8
$ws.Range("K137").Value = 7523555.399999999
$ws.Range("M137").Value = -7521005.399999999
$ws.Range("H138").Value = 5532.21
$ws.Range("I138").Value = 784.95654
$ws.Range("J138").Value = 6950.2207
$ws.Range("K138").Value = 2354.86962
$ws.Range("L138").Value = 20850.6621
$ws.Range("M138").Value = 2785.13038
$ws.Range("N138").Value = -31130.6621

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4480.375
$ws.Range("I32").Value = 3879.8823
$ws.Range("K32").Value = 3879.8823
$ws.Range("M32").Value = -3592.8823
$ws.Range("H39").Value = 36999.5
$ws.Range("I39").Value = 34000
$ws.Range("J39").Value = 39999
$ws.Range("K39").Value = 34000
$ws.Range("L39").Value = 39999
$ws.Range("M39").Value = -33480
$ws.Range("N39").Value = -41039
$ws.Range("H45").Value = 3361
$ws.Range("I45").Value = 5505.5
$ws.Range("K45").Value = 5505.5
$ws.Range("M45").Value = -5128.5
$ws.Range("H63").Value = 12596174
$ws.Range("I63").Value = 27704382
$ws.Range("J63").Value = 5999.8335
$ws.Range("K63").Value = 27704382
$ws.Range("L63").Value = 5999.8335
$ws.Range("M63").Value = -27703696
$ws.Range("N63").Value = -7371.8335
$ws.Range("H66").Value = 12596174
$ws.Range("I66").Value = 27704382
$ws.Range("J66").Value = 5999.8335
$ws.Range("K66").Value = 138521910
$ws.Range("L66").Value = 29999.1675
$ws.Range("M66").Value = -138518478
$ws.Range("N66").Value = -36863.1675
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H110").Value = 1502.381
$ws.Range("J110").Value = 1546.3334
$ws.Range("L110").Value = 1546.3334
$ws.Range("N110").Value = -5636.3334

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
$ws.Range("H94").Value = 717.92
$ws.Range("I94").Value = 630
$ws.Range("J94").Value = 1362.6666
$ws.Range("K94").Value = 630
$ws.Range("L94").Value = 1362.6666
$ws.Range("M94").Value = -179
$ws.Range("N94").Value = -2264.6666
$ws.Range("H134").Value = 2848.4666
$ws.Range("I134").Value = 2124.5454
$ws.Range("J134").Value = 4839.25
$ws.Range("K134").Value = 6373.6362
$ws.Range("L134").Value = 14517.75
$ws.Range("M134").Value = -3838.6362
$ws.Range("N134").Value = -19587.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 44195
$ws.Range("J52").Value = 44195
$ws.Range("L52").Value = 44195
$ws.Range("N52").Value = -44783
$ws.Range("H138").Value = 43072
$ws.Range("J138").Value = 43072
$ws.Range("L138").Value = 43072
$ws.Range("N138").Value = -53352
$ws.Range("H139").Value = 36666.668
$ws.Range("J139").Value = 36666.668
$ws.Range("L139").Value = 36666.668
$ws.Range("N139").Value = -46946.668

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 834905.25
$ws.Range("I5").Value = 566.8570999999999
$ws.Range("J5").Value = 1483835.1
$ws.Range("K5").Value = 1700.5713
$ws.Range("L5").Value = 4451505.300000001
$ws.Range("M5").Value = -1588.5713
$ws.Range("N5").Value = -4451729.300000001
$ws.Range("H113").Value = 3677411
$ws.Range("I113").Value = 908.9524
$ws.Range("K113").Value = 2726.8572
$ws.Range("M113").Value = -556.8571999999999
$ws.Range("H135").Value = 834905.25
$ws.Range("I135").Value = 566.8570999999999
$ws.Range("J135").Value = 1483835.1
$ws.Range("K135").Value = 5101.7139
$ws.Range("L135").Value = 13354515.9
$ws.Range("M135").Value = -2566.7139
$ws.Range("N135").Value = -13359585.9
$ws.Range("H137").Value = 2212.2632
$ws.Range("I137").Value = 2026.6666
$ws.Range("J137").Value = 2908.25
$ws.Range("K137").Value = 6079.9998
$ws.Range("L137").Value = 8724.75
$ws.Range("M137").Value = -979.9997999999996
$ws.Range("N137").Value = -18924.75

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3272.2856
$ws.Range("I132").Value = 998
$ws.Range("K132").Value = 2994
$ws.Range("M132").Value = -464

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 30400
$ws.Range("I5").Value = 30400
$ws.Range("K5").Value = 30400
$ws.Range("M5").Value = -30287
$ws.Range("H132").Value = 3883.2
$ws.Range("I132").Value = 1586
$ws.Range("J132").Value = 7770.769
$ws.Range("K132").Value = 4758
$ws.Range("L132").Value = 23312.307
$ws.Range("M132").Value = -2228
$ws.Range("N132").Value = -28372.307

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 3685.3333
$ws.Range("I38").Value = 3028
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 3028
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = -2555
$ws.Range("N38").Value = -5946
$ws.Range("H49").Value = 100000000
$ws.Range("I49").Value = 100000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 100000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -99999770
$ws.Range("N49").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 489.4
$ws.Range("I100").Value = 481.33334
$ws.Range("K100").Value = 962.66668
$ws.Range("M100").Value = -421.66668
$ws.Range("H132").Value = 22225694
$ws.Range("I132").Value = 2297.3333
$ws.Range("J132").Value = 37041292
$ws.Range("K132").Value = 6891.999899999999
$ws.Range("L132").Value = 111123876
$ws.Range("M132").Value = -4361.999899999999
$ws.Range("N132").Value = -111128936
$ws.Range("H141").Value = 36966.734
$ws.Range("J141").Value = 36966.734
$ws.Range("L141").Value = 36966.734
$ws.Range("N141").Value = -47326.734

Write-Output "Chocobo_Profits sheets refreshed."
